$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper-style inline logic: each worksheet gets one new row of log data
# appended right after its current last row. Values/formats are copied from
# the preceding row's cells where appropriate (date format, base style) so
# that the new row matches the established column typing for that sheet.
# ---------------------------------------------------------------------------

# --- Sheet "ROW50-FE-LIFTER": append row 40 ---------------------------------
$ws = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$r = 40
$ws.Cells.Item($r, 1).Value = 45744.1704925
$ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x6e"
$ws.Cells.Item($r, 5).Value = "0xe"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = 568631262647114000000000.0
$ws.Cells.Item($r, 8).Value = 366
$ws.Cells.Item($r, 9).Value = 14

# --- Sheet "ROW50-MID-LIFTER": append row 42 --------------------------------
$ws = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$r = 42
$ws.Cells.Item($r, 1).Value = 45744.13932870371
$ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat
$ws.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x72"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
# G column on this sheet stores the big integer as text (it doesn't fit
# exactly in a double), so force text formatting before assigning it,
# then restore the default style so no stray number format lingers.
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = "568631262647113771663628"
$ws.Cells.Item($r, 7).Style = "Normal"
$ws.Cells.Item($r, 8).Value = 370
$ws.Cells.Item($r, 9).Value = 25

# --- Sheet "ROW11-FE-LIFTER": append row 40 ---------------------------------
$ws = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$r = 40
$ws.Cells.Item($r, 1).Value = 45744.18933636574
$ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x6e"
$ws.Cells.Item($r, 5).Value = "0x14"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = 568631262647114000000000.0
$ws.Cells.Item($r, 8).Value = 366
$ws.Cells.Item($r, 9).Value = 20

# --- Sheet "ROW11-MID-LIFTER": append row 40 --------------------------------
$ws = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$r = 40
$ws.Cells.Item($r, 1).Value = 45744.33584197917
$ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat
$ws.Cells.Item($r, 2).Value = "0x01,0x90"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x01,0x76"
$ws.Cells.Item($r, 5).Value = "0x19"
$ws.Cells.Item($r, 6).Value = 400
$ws.Cells.Item($r, 7).Value = 568631262647114000000000.0
$ws.Cells.Item($r, 8).Value = 374
$ws.Cells.Item($r, 9).Value = 25
